# Generate Report for Handback
# The handback/regeneration run found that the de-de translation of
# "4d375c79-a9db-48a3-b453-a6097fcdfc18.md" is no longer in sync with the
# en-US source, and records new handback timestamps for it. The Overview
# sheet (and the per-language detail sheets) are updated to reflect the
# new status text, and the Status column is widened so the longer status
# string still fits.

$wb = $excel.ActiveWorkbook

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Handed back: not in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn detail sheet: row 2 is the 4d375c79-... file
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("K2").Value = "2016-09-05 10:55:48"

# ---------------------------------------------------------------------
# de-de detail sheet: row 2 is the 4d375c79-... file
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("K2").Value = "2016-09-05 10:56:09"

# ---------------------------------------------------------------------
# Overview sheet: reflect the new status for the 4d375c79-... row (row 2)
# and keep the c09d2861-... row (row 3) reporting "in sync"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $oldStatus
$wsOverview.Range("F3").Value = $oldStatus

# ---------------------------------------------------------------------
# Widen the Status columns so the longer "not in sync" text still fits
# (matches the width the workbook was regenerated with).
# ---------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 32.65
$wsOverview.Columns.Item(6).ColumnWidth = 32.65
$wsZh.Columns.Item(3).ColumnWidth = 32.65
$wsDe.Columns.Item(3).ColumnWidth = 32.65
